$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.184
$ws.Cells.Item(2, 3).Value = 0.568
$ws.Cells.Item(2, 10).Value = 0.012
$ws.Cells.Item(2, 16).Value = 0.14
$ws.Cells.Item(2, 19).Value = 0.096
$ws.Cells.Item(3, 2).Value = 0.00684931506849315
$ws.Cells.Item(3, 3).Value = 0.03424657534246575
$ws.Cells.Item(3, 10).Value = 0.03424657534246575
$ws.Cells.Item(3, 16).Value = 0.8287671232876712
$ws.Cells.Item(3, 19).Value = 0.0958904109589041
$ws.Cells.Item(4, 10).Value = 0.1290322580645161
$ws.Cells.Item(4, 16).Value = 0.5806451612903226
$ws.Cells.Item(4, 19).Value = 0.2903225806451613
$ws.Cells.Item(5, 16).Value = 0.5
$ws.Cells.Item(5, 19).Value = 0.5
$ws.Cells.Item(6, 2).Value = 0.065
$ws.Cells.Item(6, 4).Value = 0.015
$ws.Cells.Item(6, 6).Value = 0.1
$ws.Cells.Item(6, 10).Value = 0.16
$ws.Cells.Item(6, 15).Value = 0.02
$ws.Cells.Item(6, 17).Value = 0.195
$ws.Cells.Item(6, 18).Value = 0.065
$ws.Cells.Item(6, 19).Value = 0.38
$ws.Cells.Item(7, 2).Value = 0.1139240506329114
$ws.Cells.Item(7, 4).Value = 0.01265822784810127
$ws.Cells.Item(7, 5).Value = 0.006329113924050633
$ws.Cells.Item(7, 6).Value = 0.04430379746835443
$ws.Cells.Item(7, 10).Value = 0.1265822784810127
$ws.Cells.Item(7, 15).Value = 0.006329113924050633
$ws.Cells.Item(7, 17).Value = 0.1518987341772152
$ws.Cells.Item(7, 18).Value = 0.1708860759493671
$ws.Cells.Item(7, 19).Value = 0.3670886075949367
$ws.Cells.Item(8, 2).Value = 0.09711286089238845
$ws.Cells.Item(8, 4).Value = 0.007874015748031496
$ws.Cells.Item(8, 5).Value = 0.002624671916010499
$ws.Cells.Item(8, 6).Value = 0.07349081364829396
$ws.Cells.Item(8, 10).Value = 0.06036745406824147
$ws.Cells.Item(8, 15).Value = 0.01574803149606299
$ws.Cells.Item(8, 17).Value = 0.1706036745406824
$ws.Cells.Item(8, 18).Value = 0.1601049868766404
$ws.Cells.Item(8, 19).Value = 0.4120734908136483
$ws.Cells.Item(9, 2).Value = 0.07079646017699115
$ws.Cells.Item(9, 4).Value = 0.01327433628318584
$ws.Cells.Item(9, 5).Value = 0.01327433628318584
$ws.Cells.Item(9, 6).Value = 0.0752212389380531
$ws.Cells.Item(9, 10).Value = 0.1238938053097345
$ws.Cells.Item(9, 15).Value = 0.02212389380530973
$ws.Cells.Item(9, 17).Value = 0.1460176991150443
$ws.Cells.Item(9, 18).Value = 0.1194690265486726
$ws.Cells.Item(9, 19).Value = 0.415929203539823
$ws.Cells.Item(10, 2).Value = 0.09856781802864364
$ws.Cells.Item(10, 4).Value = 0.01684919966301601
$ws.Cells.Item(10, 5).Value = 0.001684919966301601
$ws.Cells.Item(10, 6).Value = 0.06571187868576242
$ws.Cells.Item(10, 10).Value = 0.109519797809604
$ws.Cells.Item(10, 15).Value = 0.01263689974726201
$ws.Cells.Item(10, 17).Value = 0.2030328559393429
$ws.Cells.Item(10, 18).Value = 0.1213142375737152
$ws.Cells.Item(10, 19).Value = 0.3706823925863522
$ws.Cells.Item(11, 7).Value = 0.1769230769230769
$ws.Cells.Item(11, 10).Value = 0.1038461538461539
$ws.Cells.Item(11, 11).Value = 0.2423076923076923
$ws.Cells.Item(11, 12).Value = 0.4576923076923077
$ws.Cells.Item(11, 19).Value = 0.01923076923076923
$ws.Cells.Item(12, 7).Value = 0.7024793388429752
$ws.Cells.Item(12, 10).Value = 0.2479338842975207
$ws.Cells.Item(12, 19).Value = 0.04958677685950413
$ws.Cells.Item(13, 7).Value = 0.6538461538461539
$ws.Cells.Item(13, 10).Value = 0.3076923076923077
$ws.Cells.Item(13, 19).Value = 0.03846153846153846
$ws.Cells.Item(15, 6).Value = 0.01401869158878505
$ws.Cells.Item(15, 8).Value = 0.2196261682242991
$ws.Cells.Item(15, 9).Value = 0.0514018691588785
$ws.Cells.Item(15, 10).Value = 0.3551401869158878
$ws.Cells.Item(15, 11).Value = 0.08411214953271028
$ws.Cells.Item(15, 13).Value = 0.009345794392523364
$ws.Cells.Item(15, 15).Value = 0.102803738317757
$ws.Cells.Item(15, 19).Value = 0.1635514018691589
$ws.Cells.Item(16, 8).Value = 0.1566265060240964
$ws.Cells.Item(16, 9).Value = 0.1325301204819277
$ws.Cells.Item(16, 10).Value = 0.3554216867469879
$ws.Cells.Item(16, 11).Value = 0.0963855421686747
$ws.Cells.Item(16, 13).Value = 0.04216867469879518
$ws.Cells.Item(16, 15).Value = 0.1024096385542169
$ws.Cells.Item(16, 19).Value = 0.1144578313253012
$ws.Cells.Item(17, 6).Value = 0.01259445843828715
$ws.Cells.Item(17, 8).Value = 0.1511335012594459
$ws.Cells.Item(17, 9).Value = 0.1360201511335013
$ws.Cells.Item(17, 10).Value = 0.4458438287153653
$ws.Cells.Item(17, 11).Value = 0.07304785894206549
$ws.Cells.Item(17, 13).Value = 0.01007556675062972
$ws.Cells.Item(17, 14).Value = 0.002518891687657431
$ws.Cells.Item(17, 15).Value = 0.04785894206549118
$ws.Cells.Item(17, 19).Value = 0.1209068010075567
$ws.Cells.Item(18, 6).Value = 0.01851851851851852
$ws.Cells.Item(18, 8).Value = 0.1592592592592593
$ws.Cells.Item(18, 9).Value = 0.09259259259259259
$ws.Cells.Item(18, 10).Value = 0.4259259259259259
$ws.Cells.Item(18, 11).Value = 0.09259259259259259
$ws.Cells.Item(18, 13).Value = 0.01851851851851852
$ws.Cells.Item(18, 14).Value = 0.003703703703703704
$ws.Cells.Item(18, 15).Value = 0.06666666666666667
$ws.Cells.Item(18, 19).Value = 0.1222222222222222
$ws.Cells.Item(19, 6).Value = 0.01216333622936577
$ws.Cells.Item(19, 8).Value = 0.1798436142484796
$ws.Cells.Item(19, 9).Value = 0.0999131190269331
$ws.Cells.Item(19, 10).Value = 0.3970460469157255
$ws.Cells.Item(19, 11).Value = 0.0947002606429192
$ws.Cells.Item(19, 13).Value = 0.02780191138140747
$ws.Cells.Item(19, 15).Value = 0.07037358818418767
$ws.Cells.Item(19, 19).Value = 0.1181581233709817
